$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.049.24"
$ws.Range("E2").Value = "  +1.35%  "

# Row 3
$ws.Range("D3").Value = "1.956.85"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'245.32"
$ws.Range("E5").Value = "  -1.14%  "

# Row 6
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").Value = "'0.4863"
$ws.Range("E7").Value = "  +0.58%  "

# Row 8
$ws.Range("D8").Value = "'0.2953"
$ws.Range("E8").Value = "  +0.50%  "

# Row 9
$ws.Range("D9").Value = "'0.06829"
$ws.Range("E9").Value = "  +0.78%  "

# Row 10
$ws.Range("D10").Value = "'19.23"
$ws.Range("E10").Value = "  -1.02%  "

# Row 11
$ws.Range("D11").Value = "'107.38"
$ws.Range("E11").Value = "  -2.74%  "

# Row 12
$ws.Range("D12").Value = "1.960.78"
$ws.Range("E12").Value = "  +0.06%  "

# Row 13
$ws.Range("D13").Value = "'0.07806"
$ws.Range("E13").Value = "  +1.01%  "

# Row 14
$ws.Range("D14").Value = "'5.460"
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("D15").Value = "'0.7039"
$ws.Range("E15").Value = "  +2.51%  "

# Row 16
$ws.Range("D16").Value = "'284.55"
$ws.Range("E16").Value = "  -2.65%  "

# Row 17
$ws.Range("D17").Value = "31.056.43"
$ws.Range("E17").Value = "  +1.33%  "

# Row 18
$ws.Range("D18").Value = "'13.20"
$ws.Range("E18").Value = "  -0.41%  "

# Row 19
$ws.Range("D19").Value = "'0.000007694"
$ws.Range("E19").Value = "  +0.42%  "

# Row 20
$ws.Range("D20").Value = "2.211.60"
$ws.Range("E20").Value = "  -0.24%  "

# Row 21
$ws.Range("E21").Value = "  +0.04%  "

# Row 22
$ws.Range("D22").Value = "'5.503"
$ws.Range("E22").Value = "  -2.43%  "

# Row 23
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("D24").Value = "'6.499"
$ws.Range("E24").Value = "  -1.42%  "

# Row 25
$ws.Range("D25").Value = "'9.811"

# Row 26
$ws.Range("D26").Value = "'170.43"
$ws.Range("E26").Value = "  -0.38%  "

# Row 27
$ws.Range("D27").Value = "'19.97"
$ws.Range("E27").Value = "  -0.79%  "

# Row 28
$ws.Range("D28").Value = "'2.200"
$ws.Range("E28").Value = "  +0.54%  "

# Row 29
$ws.Range("D29").Value = "'0.1056"
$ws.Range("E29").Value = "  -1.48%  "

# Row 30
$ws.Range("D30").Value = "'1.406"
$ws.Range("E30").Value = "  -2.42%  "

# Row 31
$ws.Range("D31").Value = "'1.584"
$ws.Range("E31").Value = "  -0.97%  "

# Row 32
$ws.Range("D32").Value = "'4.613"
$ws.Range("E32").Value = "  -1.52%  "

# Row 33
$ws.Range("D33").Value = "'4.453"
$ws.Range("E33").Value = "  +0.54%  "

# Row 34
$ws.Range("D34").Value = "'0.04931"
$ws.Range("E34").Value = "  -3.38%  "

# Row 35
$ws.Range("D35").Value = "'0.7633"
$ws.Range("E35").Value = "  -1.58%  "

# Row 36
$ws.Range("D36").Value = "'1.172"
$ws.Range("E36").Value = "  +0.22%  "

# Row 37
$ws.Range("D37").Value = "'2.732"
$ws.Range("E37").Value = "  -0.14%  "

# Row 38
$ws.Range("D38").Value = "'0.02010"
$ws.Range("E38").Value = "  -2.39%  "

# Row 39
$ws.Range("D39").Value = "'2.706"
$ws.Range("E39").Value = "  -0.20%  "

# Row 40
$ws.Range("D40").Value = "'6.544"
$ws.Range("E40").Value = "  +5.99%  "

# Row 41
$ws.Range("D41").Value = "'2.103"
$ws.Range("E41").Value = "  +1.60%  "

# Row 42
$ws.Range("D42").Value = "'75.06"
$ws.Range("E42").Value = "  +7.37%  "

# Row 43
$ws.Range("D43").Value = "'0.8903"
$ws.Range("E43").Value = "  +1.79%  "

# Row 44
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.4464"
$ws.Range("E44").Value = "  +0.03%  "

# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'109.37"
$ws.Range("E45").Value = "  -0.48%  "

# Row 46
$ws.Range("D46").Value = "'8.199"
$ws.Range("E46").Value = "  +11.02%  "

# Row 47
$ws.Range("E47").Value = "  -0.03%  "

# Row 48
$ws.Range("D48").Value = "1.003.59"
$ws.Range("E48").Value = "  +11.17%  "

# Row 49
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1255"
$ws.Range("E49").Value = "  -1.74%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.299"
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("D51").Value = "'35.76"
$ws.Range("E51").Value = "  -0.47%  "
